$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point (2026/01/22, 木曜日, 13:00, ranking 154) was recorded and
# inserted into the time-ordered log at row 703, pushing the previously
# existing rows 703-744 down to 704-745 (dimension grows from D744 to D745).
$ws.Rows.Item(703).Insert()

# Write the new row. Column A holds a date-like string that must stay as
# literal text (matching every other date cell in the sheet, which are plain
# text, not real Excel dates) -- force Text format before assigning the
# value so it isn't auto-parsed into a date serial, then clear the
# number-format override back off so the cell ends up with the sheet's
# default (unstyled) look, same as its neighbours.
$ws.Range("A703").NumberFormat = "@"
$ws.Range("A703").Value = "2026/01/22"
$ws.Range("A703").ClearFormats()

$ws.Range("B703").Value = "木"
$ws.Range("C703").Value = 13
$ws.Range("D703").Value = 154
